$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert first new record at row 27 ---
# This pushes the existing rows 27-37 down to 28-38.
$ws.Rows("27").Insert()

$ws.Range("A27").Value = 11
$ws.Range("B27").Value = "Vega Monumental Concepción"
$ws.Range("C27").Value = "Bíobío"
$ws.Range("D27").Value = 44435
$ws.Range("E27").Value = 8
$ws.Range("F27").Value = 100112021
$ws.Range("G27").Value = "Ají"
$ws.Range("H27").Value = "Inferno"
$ws.Range("I27").Value = "Primera"
$ws.Range("J27").Value = 40
$ws.Range("K27").Value = 37000
$ws.Range("L27").Value = 38000
$ws.Range("M27").Value = 37500
$ws.Range("N27").Value = "`$/caja 12 kilos"
$ws.Range("O27").Value = "Región de Arica y Parinacota"
$ws.Range("P27").Value = 3125
$ws.Range("Q27").Value = 12
$ws.Range("R27").Value = "Hortaliza"

# --- Insert second new record at row 32 ---
# After the first insert, the old row 31 now sits at row 32, so inserting
# here pushes the current rows 32-38 down to 33-39.
$ws.Rows("32").Insert()

$ws.Range("A32").Value = 11
$ws.Range("B32").Value = "Vega Monumental Concepción"
$ws.Range("C32").Value = "Bíobío"
$ws.Range("D32").Value = 44433
$ws.Range("E32").Value = 8
$ws.Range("F32").Value = 100112021
$ws.Range("G32").Value = "Ají"
$ws.Range("H32").Value = "Inferno"
$ws.Range("I32").Value = "Primera"
$ws.Range("J32").Value = 40
$ws.Range("K32").Value = 37000
$ws.Range("L32").Value = 38000
$ws.Range("M32").Value = 37500
$ws.Range("N32").Value = "`$/caja 12 kilos"
$ws.Range("O32").Value = "Región de Arica y Parinacota"
$ws.Range("P32").Value = 3125
$ws.Range("Q32").Value = 12
$ws.Range("R32").Value = "Hortaliza"
